$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ncam1", "Robo3", "ECs", 3, 1, 0.9949870000000001, 2.984961, 0.03855738270564991, 0.03855738270564991, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 1.320145434976667, 11.88130891479, 0.03841902983729085, 0.03841902983729086),
    @("ECs", "Ncam1", "Robo3", "MuSCs", 3, 1, 0.9949870000000001, 2.984961, 0.03855738270564991, 0.03855738270564991, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.004754047886000001, 0.042786430974, 0.0001383528683590621, 0.0001383528683590621),
    @("FAPs", "Ncam1", "Robo3", "ECs", 3, 1, 8.469728666666667, 25.409186, 0.3282159160005915, 0.3282159160005916, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 11.23760776250444, 101.13846986254, 0.3270382008593321, 0.3270382008593322),
    @("FAPs", "Ncam1", "Robo3", "MuSCs", 3, 1, 8.469728666666667, 25.409186, 0.3282159160005915, 0.3282159160005916, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.04046836356933334, 0.364215272124, 0.001177715141259441, 0.001177715141259441),
    @("Inflammatory-Mac", "Ncam1", "Robo3", "ECs", 2, 0.6666666666666666, 0.04495399999999999, 0.134862, 0.001742041435867791, 0.001742041435867791, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 0.05964481735333332, 0.5368033561799999, 0.001735790585510738, 0.001735790585510738),
    @("Inflammatory-Mac", "Ncam1", "Robo3", "MuSCs", 2, 0.6666666666666666, 0.04495399999999999, 0.134862, 0.001742041435867791, 0.001742041435867791, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.000214790212, 0.001933111908, 0.000006250850357053183, 0.000006250850357053183),
    @("MuSCs", "Ncam1", "Robo3", "ECs", 3, 1, 15.96019966666667, 47.880599, 0.6184839868322428, 0.6184839868322429, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 21.17593971706778, 190.58345745361, 0.6162647222554527, 0.6162647222554528),
    @("MuSCs", "Ncam1", "Robo3", "MuSCs", 3, 1, 15.96019966666667, 47.880599, 0.6184839868322428, 0.6184839868322429, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.07625783400733334, 0.686320506066, 0.002219264576790129, 0.002219264576790129),
    @("Neutrophils", "Ncam1", "Robo3", "ECs", 3, 1, 0.3268106666666666, 0.980432, 0.01266445083901121, 0.01266445083901121, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 0.4336113031644444, 3.90250172848, 0.01261900784011407, 0.01261900784011407),
    @("Neutrophils", "Ncam1", "Robo3", "MuSCs", 3, 1, 0.3268106666666666, 0.980432, 0.01266445083901121, 0.01266445083901121, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.001561501365333333, 0.014053512288, 0.00004544299889714201, 0.00004544299889714201),
    @("Resolving-Mac", "Ncam1", "Robo3", "ECs", 1, 0.3333333333333333, 0.008676333333333333, 0.026029, 0.0003362221866367304, 0.0003362221866367304, 3, 1, 1.326796666666667, 3.98039, 0.9964117671208325, 0.9964117671208325, 0.01151173014555555, 0.10360557131, 0.0003350157431319349, 0.0003350157431319349),
    @("Resolving-Mac", "Ncam1", "Robo3", "MuSCs", 1, 0.3333333333333333, 0.008676333333333333, 0.026029, 0.0003362221866367304, 0.0003362221866367304, 1, 0.3333333333333333, 0.004778, 0.014334, 0.003588232879167623, 0.003588232879167622, 0.00004145552066666666, 0.000373099686, 0.000001206443504795549, 0.000001206443504795549),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

